$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, value first
$ws.Range("H1").Value = "Save"

# Clone the header formatting (bold, centered, bordered) from the existing
# "sum" header (G1) onto the new "Save" header (H1) using copy/paste-special
# so the engine dedupes onto the same style index rather than minting a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add values for the new Save column (row 2 -> 1, row 3 -> 0)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
